$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H28").Value = 761.3333
$ws_ALC.Range("I28").Value = 317.27274
$ws_ALC.Range("J28").Value = 1249.8
$ws_ALC.Range("K28").Value = 317.27274
$ws_ALC.Range("L28").Value = 1249.8
$ws_ALC.Range("M28").Value = 167.72726
$ws_ALC.Range("N28").Value = -2219.8
$ws_ALC.Range("H137").Value = 1733.0541
$ws_ALC.Range("I137").Value = 1692.92
$ws_ALC.Range("J137").Value = 1816.6666
$ws_ALC.Range("K137").Value = 5078.76
$ws_ALC.Range("L137").Value = 5449.9998
$ws_ALC.Range("M137").Value = -2528.76
$ws_ALC.Range("N137").Value = -10549.9998
$ws_ALC.Range("H138").Value = 13516142
$ws_ALC.Range("I138").Value = 33334560
$ws_ALC.Range("J138").Value = 3584.0908
$ws_ALC.Range("K138").Value = 100003680
$ws_ALC.Range("L138").Value = 10752.2724
$ws_ALC.Range("M138").Value = -99998540
$ws_ALC.Range("N138").Value = -21032.2724
$ws_ALC.Range("H141").Value = 1187.3448
$ws_ALC.Range("I141").Value = 703.2041
$ws_ALC.Range("J141").Value = 3823.2222
$ws_ALC.Range("K141").Value = 2109.6123
$ws_ALC.Range("L141").Value = 11469.6666
$ws_ALC.Range("M141").Value = 3070.3877
$ws_ALC.Range("N141").Value = -21829.6666
$ws_ARM.Range("H2").Value = 2256.25
$ws_ARM.Range("I2").Value = 1864.2858
$ws_ARM.Range("K2").Value = 1864.2858
$ws_ARM.Range("M2").Value = -1751.2858
$ws_ARM.Range("H32").Value = 3899.8572
$ws_ARM.Range("I32").Value = 2548.75
$ws_ARM.Range("J32").Value = 20113.143
$ws_ARM.Range("K32").Value = 2548.75
$ws_ARM.Range("L32").Value = 20113.143
$ws_ARM.Range("M32").Value = -2261.75
$ws_ARM.Range("N32").Value = -20687.143
$ws_ARM.Range("H61").Value = 340722.47
$ws_ARM.Range("I61").Value = 368275.28
$ws_ARM.Range("J61").Value = 3200.5
$ws_ARM.Range("K61").Value = 368275.28
$ws_ARM.Range("L61").Value = 3200.5
$ws_ARM.Range("M61").Value = -368063.28
$ws_ARM.Range("N61").Value = -3624.5
$ws_ARM.Range("H116").Value = 2256.25
$ws_ARM.Range("I116").Value = 1864.2858
$ws_ARM.Range("K116").Value = 1864.2858
$ws_ARM.Range("M116").Value = 429.7141999999999
$ws_ARM.Range("H132").Value = 11410.86
$ws_ARM.Range("I132").Value = 1393.878
$ws_ARM.Range("J132").Value = 57043.777
$ws_ARM.Range("K132").Value = 4181.634
$ws_ARM.Range("L132").Value = 171131.331
$ws_ARM.Range("M132").Value = -1651.634
$ws_ARM.Range("N132").Value = -176191.331
$ws_ARM.Range("H136").Value = 340722.47
$ws_ARM.Range("I136").Value = 368275.28
$ws_ARM.Range("J136").Value = 3200.5
$ws_ARM.Range("K136").Value = 1104825.84
$ws_ARM.Range("L136").Value = 9601.5
$ws_ARM.Range("M136").Value = -1102275.84
$ws_ARM.Range("N136").Value = -14701.5
$ws_BSM.Range("H3").Value = 2256.25
$ws_BSM.Range("I3").Value = 1864.2858
$ws_BSM.Range("K3").Value = 1864.2858
$ws_BSM.Range("M3").Value = -1750.2858
$ws_BSM.Range("H105").Value = 11400
$ws_BSM.Range("I105").Value = 19000
$ws_BSM.Range("J105").Value = 3800
$ws_BSM.Range("K105").Value = 19000
$ws_BSM.Range("L105").Value = 3800
$ws_BSM.Range("M105").Value = -17253
$ws_BSM.Range("N105").Value = -7294
$ws_BSM.Range("H107").Value = 819.63635
$ws_BSM.Range("J107").Value = 1571
$ws_BSM.Range("L107").Value = 1571
$ws_BSM.Range("N107").Value = -5411
$ws_BSM.Range("H134").Value = 3067.628
$ws_BSM.Range("I134").Value = 3815.4482
$ws_BSM.Range("K134").Value = 11446.3446
$ws_BSM.Range("M134").Value = -8911.3446
$ws_CRP.Range("H16").Value = 848.75
$ws_CRP.Range("I16").Value = 761.4545000000001
$ws_CRP.Range("J16").Value = 1040.8
$ws_CRP.Range("K16").Value = 761.4545000000001
$ws_CRP.Range("L16").Value = 1040.8
$ws_CRP.Range("M16").Value = -474.4545000000001
$ws_CRP.Range("N16").Value = -1614.8
$ws_CRP.Range("H31").Value = 4778.8335
$ws_CRP.Range("I31").Value = 2942.077
$ws_CRP.Range("K31").Value = 2942.077
$ws_CRP.Range("M31").Value = -2647.077
$ws_CRP.Range("H34").Value = 4778.8335
$ws_CRP.Range("I34").Value = 2942.077
$ws_CRP.Range("K34").Value = 2942.077
$ws_CRP.Range("M34").Value = -2740.077
$ws_CRP.Range("H58").Value = 821.65753
$ws_CRP.Range("I58").Value = 658.24
$ws_CRP.Range("J58").Value = 1176.9131
$ws_CRP.Range("K58").Value = 658.24
$ws_CRP.Range("L58").Value = 1176.9131
$ws_CRP.Range("M58").Value = -455.24
$ws_CRP.Range("N58").Value = -1582.9131
$ws_CRP.Range("H105").Value = 11364630
$ws_CRP.Range("I105").Value = 13889736
$ws_CRP.Range("K105").Value = 13889736
$ws_CRP.Range("M105").Value = -13887989
$ws_CRP.Range("H113").Value = 848.75
$ws_CRP.Range("I113").Value = 761.4545000000001
$ws_CRP.Range("J113").Value = 1040.8
$ws_CRP.Range("K113").Value = 761.4545000000001
$ws_CRP.Range("L113").Value = 1040.8
$ws_CRP.Range("M113").Value = 1408.5455
$ws_CRP.Range("N113").Value = -5380.8
$ws_CRP.Range("H132").Value = 2562.3872
$ws_CRP.Range("I132").Value = 1947.3334
$ws_CRP.Range("K132").Value = 5842.0002
$ws_CRP.Range("M132").Value = -3312.0002
$ws_CRP.Range("H136").Value = 821.65753
$ws_CRP.Range("I136").Value = 658.24
$ws_CRP.Range("J136").Value = 1176.9131
$ws_CRP.Range("K136").Value = 1974.72
$ws_CRP.Range("L136").Value = 3530.7393
$ws_CRP.Range("M136").Value = 575.28
$ws_CRP.Range("N136").Value = -8630.739300000001
$ws_CUL.Range("H5").Value = 1749.2
$ws_CUL.Range("I5").Value = 1532.3334
$ws_CUL.Range("J5").Value = 1893.7778
$ws_CUL.Range("K5").Value = 4597.0002
$ws_CUL.Range("L5").Value = 5681.3334
$ws_CUL.Range("M5").Value = -4485.0002
$ws_CUL.Range("N5").Value = -5905.3334
$ws_CUL.Range("H34").Value = 1265.3334
$ws_CUL.Range("I34").Value = 414
$ws_CUL.Range("J34").Value = 1549.1111
$ws_CUL.Range("K34").Value = 1242
$ws_CUL.Range("L34").Value = 4647.3333
$ws_CUL.Range("M34").Value = -1158
$ws_CUL.Range("N34").Value = -4815.3333
$ws_CUL.Range("H39").Value = 2940
$ws_CUL.Range("J39").Value = 2940
$ws_CUL.Range("L39").Value = 8820
$ws_CUL.Range("N39").Value = -9408
$ws_CUL.Range("H118").Value = 55557348
$ws_CUL.Range("I118").Value = 83333560
$ws_CUL.Range("J118").Value = 4932.6665
$ws_CUL.Range("K118").Value = 250000680
$ws_CUL.Range("L118").Value = 14797.9995
$ws_CUL.Range("M118").Value = -249999437
$ws_CUL.Range("N118").Value = -17283.9995
$ws_CUL.Range("H131").Value = 213678.75
$ws_CUL.Range("J131").Value = 250953.53
$ws_CUL.Range("L131").Value = 752860.59
$ws_CUL.Range("N131").Value = -762940.59
$ws_CUL.Range("H135").Value = 1749.2
$ws_CUL.Range("I135").Value = 1532.3334
$ws_CUL.Range("J135").Value = 1893.7778
$ws_CUL.Range("K135").Value = 13791.0006
$ws_CUL.Range("L135").Value = 17044.0002
$ws_CUL.Range("M135").Value = -11256.0006
$ws_CUL.Range("N135").Value = -22114.0002
$ws_CUL.Range("H138").Value = 1722.2354
$ws_CUL.Range("I138").Value = 1448.5714
$ws_CUL.Range("J138").Value = 2999.3333
$ws_CUL.Range("K138").Value = 4345.7142
$ws_CUL.Range("L138").Value = 8997.999899999999
$ws_CUL.Range("M138").Value = 794.2857999999997
$ws_CUL.Range("N138").Value = -19277.9999
$ws_GSM.Range("H113").Value = 4260
$ws_GSM.Range("I113").Value = 2650
$ws_GSM.Range("J113").Value = 5333.3335
$ws_GSM.Range("K113").Value = 2650
$ws_GSM.Range("L113").Value = 5333.3335
$ws_GSM.Range("M113").Value = -480
$ws_GSM.Range("N113").Value = -9673.333500000001
$ws_GSM.Range("H132").Value = 15384.368
$ws_GSM.Range("I132").Value = 2414.879
$ws_GSM.Range("J132").Value = 100983
$ws_GSM.Range("K132").Value = 7244.637
$ws_GSM.Range("L132").Value = 302949
$ws_GSM.Range("M132").Value = -4714.637
$ws_GSM.Range("N132").Value = -308009
$ws_LTW.Range("H132").Value = 1137.4706
$ws_LTW.Range("I132").Value = 1102.2554
$ws_LTW.Range("J132").Value = 1551.25
$ws_LTW.Range("K132").Value = 3306.7662
$ws_LTW.Range("L132").Value = 4653.75
$ws_LTW.Range("M132").Value = -776.7662
$ws_LTW.Range("N132").Value = -9713.75
$ws_WVR.Range("H132").Value = 803.95123
$ws_WVR.Range("I132").Value = 673.2
$ws_WVR.Range("J132").Value = 1566.6666
$ws_WVR.Range("K132").Value = 2019.6
$ws_WVR.Range("L132").Value = 4699.9998
$ws_WVR.Range("M132").Value = 510.3999999999999
$ws_WVR.Range("N132").Value = -9759.9998
$ws_WVR.Range("H136").Value = 13334916
$ws_WVR.Range("I136").Value = 18868682
$ws_WVR.Range("J136").Value = 3568.8635
$ws_WVR.Range("K136").Value = 56606046
$ws_WVR.Range("L136").Value = 10706.5905
$ws_WVR.Range("M136").Value = -56603496
$ws_WVR.Range("N136").Value = -15806.5905
$ws_ARM.Range("H102").Value = 1192.5
$ws_ARM.Range("I102").Value = 1192.5
$ws_ARM.Range("J102").Value = 0
$ws_ARM.Range("K102").Value = 1192.5
$ws_ARM.Range("L102").Value = 0
$ws_ARM.Range("M102").Value = 429.5
$ws_ARM.Range("N102").ClearContents()
$ws_CUL.Range("H55").Value = 2642.8572
$ws_CUL.Range("I55").Value = 0
$ws_CUL.Range("J55").Value = 2642.8572
$ws_CUL.Range("K55").Value = 0
$ws_CUL.Range("L55").Value = 7928.571599999999
$ws_CUL.Range("M55").ClearContents()
$ws_CUL.Range("N55").Value = -8282.571599999999
$ws_GSM.Range("H138").Value = 50429
$ws_GSM.Range("J138").Value = 50429
$ws_GSM.Range("L138").Value = 50429
$ws_GSM.Range("N138").Value = -60709
